# Update column F (dSF) values for the rows that were re-pulled / recalculated.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = 3
    3  = 3
    4  = 1
    5  = -2
    6  = 6
    7  = -4
    8  = 1
    12 = 3
    13 = -4
    14 = 5
    15 = 0
    16 = 6
    17 = 6
    18 = 1
    19 = -2
    20 = 4
    21 = 3
    23 = -5
    24 = -3
    25 = -2
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
